# VerbTense2 Repo.xlsx - "Add files via upload" edit
# 1) Flip MatchCase ("C" column) from "N" to "Y" for a set of verbs.
# 2) Insert two brand-new verb rows: "deploying"/"deploy" and "solving"/"solve".
# 3) Refresh the sheet view (zoom + selection) to match the re-saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1) Flip MatchCase to "Y" for these verbs (using the ORIGINAL, pre-insert row numbers) ---
$yRows = @(17, 22, 32, 34, 35, 42, 47, 51, 78, 84)
foreach ($r in $yRows) {
    $ws.Cells.Item($r, 3).Value = '"Y"'
}

# --- 2) Insert "deploying" / "deploy" row right after "designing" (original row 21) ---
$ws.Rows(22).Insert()
$ws.Cells.Item(22, 1).Value = '"deploying"'
$ws.Cells.Item(22, 2).Value = '"deploy"'
$ws.Cells.Item(22, 3).Value = '"N"'
$ws.Cells.Item(22, 4).Formula = '="verbTense(" & ROW(A22)-1 & ", " & 1 & ") = " & A22'
$ws.Cells.Item(22, 5).Formula = '="verbTense(" & ROW(A22)-1 & ", " & 2 & ") = " & B22'
$ws.Cells.Item(22, 6).Formula = '="verbTense(" & ROW(A22)-1 & ", " & 3 & ") = " & C22'
$ws.Cells.Item(22, 7).Formula = '=D22 & " : " & E22 & " : " & F22'

# --- 3) Insert "solving" / "solve" row right after "sharing" (original row 74, now row 75) ---
$ws.Rows(76).Insert()
$ws.Cells.Item(76, 1).Value = '"solving"'
$ws.Cells.Item(76, 2).Value = '"solve"'
$ws.Cells.Item(76, 3).Value = '"N"'
$ws.Cells.Item(76, 4).Formula = '="verbTense(" & ROW(A76)-1 & ", " & 1 & ") = " & A76'
$ws.Cells.Item(76, 5).Formula = '="verbTense(" & ROW(A76)-1 & ", " & 2 & ") = " & B76'
$ws.Cells.Item(76, 6).Formula = '="verbTense(" & ROW(A76)-1 & ", " & 3 & ") = " & C76'
$ws.Cells.Item(76, 7).Formula = '=D76 & " : " & E76 & " : " & F76'

# --- 4) Refresh the view: zoomed to 55%, selection on column G's data range ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 55
$ws.Range("G2:G87").Select()
